# Loan RBI, Variable Instalments
# Inserts a new (blank) column before column N on the "Repayment schedule"
# sheet, shifting the existing N:P columns (Late / heading / Outstanding)
# one column to the right, to O:Q.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

$ws.Activate()

# Insert a new column at N, shifting N:P (and everything after) to the right.
$ws.Columns("N").Insert(-4161)

# The freshly inserted column inherits column M's width (matches Excel's
# own behaviour of carrying the left neighbour's width onto a new column).
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Re-stamp the "Outstanding" column (old P, now Q) with its exact original
# values -- the shift can reintroduce binary floating point noise
# (875.82 -> 875.82000000000005) even though no formula is involved.
$outstanding = @{
    3  = 0
    4  = 875.82
    5  = 805.89
    6  = 857.48
    7  = 831.17
    8  = 825.64
    9  = 816.11
    10 = 810.08
    11 = 802.3
    12 = 793.52
    13 = 786.74
    14 = 778.46
    15 = 771.22
}
foreach ($row in $outstanding.Keys) {
    $ws.Cells.Item($row, 17).Value = $outstanding[$row]
}

# Update the selection to match the authored workbook (S8 on the active sheet).
$ws.Range("S8").Select()
